$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.2969342779685747, 0.6523826063768725, 0.6631625730242209, 0.8143479434641073, 0.7790613848385668, 19)
    3  = @(0.1798162647974569, 0.6425775241922705, 0.6172442385525494, 0.7856489283086622, 0.7869668611904497, 18)
    4  = @(0.1712112111122847, 0.6344048150359224, 0.6409194659967445, 0.8005744599952864, 0.8061213146204157, 17)
    5  = @(0.3430733728802315, 0.6296856560655977, 0.5655510257999548, 0.752031266504229,  0.6911645721017333, 16)
    6  = @(0.2773126667199611, 0.591515179616282,  0.5117981096136305, 0.7154006636938706, 0.6826124667855324, 15)
    7  = @(0.3190424291774069, 0.5868986679421527, 0.5229048546182766, 0.7231216043088995, 0.6734317943900356, 14)
    8  = @(0.3116543284621223, 0.6160086454085564, 0.559373466974749,  0.7479127402142238, 0.7076478411394824, 13)
    9  = @(0.4347727426752494, 0.5707397934898671, 0.4730317857967594, 0.687773062715282,  0.5566174936231759, 12)
    10 = @(0.3216259016166307, 0.4787523040349611, 0.2822515447740024, 0.531273512208168,  0.443496512502208,  11)
    11 = @(0.4110739536262513, 0.4543498161016222, 0.2611419507433211, 0.5110204993376695, 0.3200002697801009, 10)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
